{"js": "// Word JS API (Office.js) edit script.\n// Body of: async (context) => { ... }\n//\n// Three wording tweaks inside the \"Situation-Orientation\" / \"Future-Oriented\n// Organisation\" table rows:\n//   1) Merge \"The Agile way of working ... Scrum activites. \" + \"Requirement\n//      gathering\" + \" was also a familiar activity ... In addition, a \" into\n//      a single contiguous run of text (content itself is unchanged).\n//   2) \"Version Control even though is familiar still needs more practice...\"\n//      -> \"Version Control, even though familiar, still needs more practice...\"\n//   3) \"Regarding planning and follow through it has been suceessful until\n//      this point.\" -> \"Planning and follow through have been suceessful\n//      until this point.\"\n\nconst body = context.document.body;\n\n// --- Change 1 -----------------------------------------------------------\nconst search1 = body.search(\n  \"The Agile way of working is already know and thus was easy to get onboard with Scrum activites. Requirement gathering was also a familiar activity due to various green projects done in the past. In addition, a \",\n  { matchCase: true }\n);\nsearch1.load(\"items\");\nawait context.sync();\n\nif (search1.items.length > 0) {\n  const range1 = search1.items[0];\n  range1.insertText(\n    \"The Agile way of working is already know and thus was easy to get onboard with Scrum activites. Requirement gathering was also a familiar activity due to various green projects done in the past. In addition, a \",\n    Word.InsertLocation.replace\n  );\n}\n\n// --- Change 2 -------------------------------------------------------------\nconst search2 = body.search(\n  \"Version Control even though is familiar still needs more practice in the context of a larger team and better procedures.\",\n  { matchCase: true }\n);\nsearch2.load(\"items\");\nawait context.sync();\n\nif (search2.items.length > 0) {\n  const range2 = search2.items[0];\n  range2.insertText(\n    \"Version Control, even though familiar, still needs more practice in the context of a larger team and better procedures.\",\n    Word.InsertLocation.replace\n  );\n}\n\n// --- Change 3 -------------------------------------------------------------\nconst search3 = body.search(\n  \"Regarding planning and follow through it has been suceessful until this point.\",\n  { matchCase: true }\n);\nsearch3.load(\"items\");\nawait context.sync();\n\nif (search3.items.length > 0) {\n  const range3 = search3.items[0];\n  range3.insertText(\n    \"Planning and follow through have been suceessful until this point.\",\n    Word.InsertLocation.replace\n  );\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop edit script.\n# The document is open as $word.ActiveDocument.\n#\n# Three wording tweaks inside the \"Situation-Orientation\" / \"Future-Oriented\n# Organisation\" table rows:\n#   1) Merge \"The Agile way of working ... Scrum activites. \" + \"Requirement\n#      gathering\" + \" was also a familiar activity ... In addition, a \" into\n#      a single contiguous run of text (content itself is unchanged).\n#   2) \"Version Control even though is familiar still needs more practice...\"\n#      -> \"Version Control, even though familiar, still needs more practice...\"\n#   3) \"Regarding planning and follow through it has been suceessful until\n#      this point.\" -> \"Planning and follow through have been suceessful\n#      until this point.\"\n\n$d = $word.ActiveDocument\n\n# wdReplaceAll = 2\n$wdReplaceAll = 2\n\n# --- Change 1 --------------------------------------------------------------\n$rng1 = $d.Content\n$rng1.Find.Execute(\n  \"The Agile way of working is already know and thus was easy to get onboard with Scrum activites. Requirement gathering was also a familiar activity due to various green projects done in the past. In addition, a \",\n  $false, $false, $false, $false, $false, $true, 1, $false,\n  \"The Agile way of working is already know and thus was easy to get onboard with Scrum activites. Requirement gathering was also a familiar activity due to various green projects done in the past. In addition, a \",\n  $wdReplaceAll\n)\n\n# --- Change 2 --------------------------------------------------------------\n$rng2 = $d.Content\n$rng2.Find.Execute(\n  \"Version Control even though is familiar still needs more practice in the context of a larger team and better procedures.\",\n  $false, $false, $false, $false, $false, $true, 1, $false,\n  \"Version Control, even though familiar, still needs more practice in the context of a larger team and better procedures.\",\n  $wdReplaceAll\n)\n\n# --- Change 3 --------------------------------------------------------------\n$rng3 = $d.Content\n$rng3.Find.Execute(\n  \"Regarding planning and follow through it has been suceessful until this point.\",\n  $false, $false, $false, $false, $false, $true, 1, $false,\n  \"Planning and follow through have been suceessful until this point.\",\n  $wdReplaceAll\n)\n"}
